$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.505.70"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "2.478.14"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.74%  "

$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("D9").Value = "2.477.98"
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0988"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.325"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.37%  "

$ws.Range("D14").Value = "2.925.03"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").Value = "58.384.39"
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").Value = "2.492.09"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "318.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.10%  "

$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.403"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("E26").Value = "  +1.25%  "

$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.50%  "

$ws.Range("D30").Value = "0.0₃0743"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").Value = "  +0.66%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.35%  "

$ws.Range("E39").Value = "  +2.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "270.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "128.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.585"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0933"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0498"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0215"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").Value = "1.730.21"
$ws.Range("E51").Value = "  -0.35%  "
